$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 151.8
$ws.Range("I2").Value = 151.8
$ws.Range("K2").Value = 151.8
$ws.Range("M2").Value = -38.80000000000001
$ws.Range("H9").Value = 237.93333
$ws.Range("J9").Value = 281.75
$ws.Range("L9").Value = 281.75
$ws.Range("N9").Value = -619.75
$ws.Range("H32").Value = 7719.091
$ws.Range("I32").Value = 6831
$ws.Range("J32").Value = 8052.125
$ws.Range("K32").Value = 6831
$ws.Range("L32").Value = 8052.125
$ws.Range("M32").Value = -6505
$ws.Range("N32").Value = -8704.125
$ws.Range("H41").Value = 489
$ws.Range("I41").Value = 301.5
$ws.Range("K41").Value = 301.5
$ws.Range("M41").Value = 138.5
$ws.Range("H43").Value = 3517.7778
$ws.Range("J43").Value = 1534.2
$ws.Range("L43").Value = 1534.2
$ws.Range("N43").Value = -1672.2
$ws.Range("H80").Value = 527.1667
$ws.Range("J80").Value = 805.9167
$ws.Range("L80").Value = 2417.7501
$ws.Range("N80").Value = -4413.7501
$ws.Range("H83").Value = 527.1667
$ws.Range("J83").Value = 805.9167
$ws.Range("L83").Value = 7253.2503
$ws.Range("N83").Value = -17237.2503
$ws.Range("H92").Value = 1199.7142
$ws.Range("I92").Value = 944.7
$ws.Range("K92").Value = 944.7
$ws.Range("M92").Value = 303.3
$ws.Range("H96").Value = 1500
$ws.Range("I96").Value = 1500
$ws.Range("K96").Value = 4500
$ws.Range("M96").Value = -3127
$ws.Range("H101").Value = 1149.0834
$ws.Range("I101").Value = 649.75
$ws.Range("J101").Value = 1398.75
$ws.Range("K101").Value = 1949.25
$ws.Range("L101").Value = 4196.25
$ws.Range("M101").Value = -327.25
$ws.Range("N101").Value = -7440.25
$ws.Range("H109").Value = 69995
$ws.Range("J109").Value = 69995
$ws.Range("L109").Value = 69995
$ws.Range("N109").Value = -72769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2331.4119
$ws.Range("I32").Value = 1467.5217
$ws.Range("K32").Value = 1467.5217
$ws.Range("M32").Value = -1180.5217
$ws.Range("H80").Value = 50110
$ws.Range("J80").Value = 50110
$ws.Range("L80").Value = 50110
$ws.Range("N80").Value = -52106
$ws.Range("H83").Value = 50110
$ws.Range("J83").Value = 50110
$ws.Range("L83").Value = 150330
$ws.Range("N83").Value = -160314

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 57500
$ws.Range("J62").Value = 45000
$ws.Range("L62").Value = 45000
$ws.Range("N62").Value = -46372
$ws.Range("H65").Value = 57500
$ws.Range("J65").Value = 45000
$ws.Range("L65").Value = 135000
$ws.Range("N65").Value = -141864
$ws.Range("H99").Value = 100003110
$ws.Range("I99").Value = 111114344
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 111114344
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -111112846
$ws.Range("N99").Value = -4996
$ws.Range("H105").Value = 2122.6875
$ws.Range("I105").Value = 2067.5334
$ws.Range("K105").Value = 2067.5334
$ws.Range("M105").Value = -320.5333999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1674.7333
$ws.Range("I31").Value = 2556
$ws.Range("K31").Value = 2556
$ws.Range("M31").Value = -2261
$ws.Range("H34").Value = 1674.7333
$ws.Range("I34").Value = 2556
$ws.Range("K34").Value = 2556
$ws.Range("M34").Value = -2354
$ws.Range("H82").Value = 24995
$ws.Range("I82").Value = 24995
$ws.Range("K82").Value = 24995
$ws.Range("M82").Value = -24634
$ws.Range("H85").Value = 24995
$ws.Range("I85").Value = 24995
$ws.Range("K85").Value = 24995
$ws.Range("M85").Value = -23747
$ws.Range("H87").Value = 40000
$ws.Range("J87").Value = 40000
$ws.Range("L87").Value = 40000
$ws.Range("N87").Value = -42372
$ws.Range("H88").Value = 24716.25
$ws.Range("J88").Value = 24716.25
$ws.Range("L88").Value = 24716.25
$ws.Range("N88").Value = -25528.25
$ws.Range("H90").Value = 40000
$ws.Range("J90").Value = 40000
$ws.Range("L90").Value = 120000
$ws.Range("N90").Value = -131856
$ws.Range("H91").Value = 24716.25
$ws.Range("J91").Value = 24716.25
$ws.Range("L91").Value = 24716.25
$ws.Range("N91").Value = -27524.25
$ws.Range("H107").Value = 1083.8928
$ws.Range("I107").Value = 627.5
$ws.Range("J107").Value = 1159.9584
$ws.Range("K107").Value = 627.5
$ws.Range("L107").Value = 1159.9584
$ws.Range("M107").Value = 1292.5
$ws.Range("N107").Value = -4999.9584
$ws.Range("H132").Value = 2293.25
$ws.Range("I132").Value = 2044.6451
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 6133.9353
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -3603.9353
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 62500510
$ws.Range("I4").Value = 76923480
$ws.Range("K4").Value = 230770440
$ws.Range("M4").Value = -230770328
$ws.Range("H51").Value = 4175
$ws.Range("I51").Value = 3300
$ws.Range("J51").Value = 4466.6665
$ws.Range("K51").Value = 9900
$ws.Range("L51").Value = 13399.9995
$ws.Range("M51").Value = -9440
$ws.Range("N51").Value = -14319.9995
$ws.Range("H58").Value = 5
$ws.Range("I58").Value = 5
$ws.Range("K58").Value = 15
$ws.Range("M58").Value = 113

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 21188
$ws.Range("I70").Value = 34536.277
$ws.Range("K70").Value = 34536.277
$ws.Range("M70").Value = -34266.277
$ws.Range("H73").Value = 21188
$ws.Range("I73").Value = 34536.277
$ws.Range("K73").Value = 34536.277
$ws.Range("M73").Value = -33600.277
$ws.Range("H97").Value = 674.2727
$ws.Range("I97").Value = 531.6
$ws.Range("J97").Value = 980
$ws.Range("K97").Value = 531.6
$ws.Range("L97").Value = 980
$ws.Range("M97").Value = -35.60000000000002
$ws.Range("N97").Value = -1972
$ws.Range("H102").Value = 2724.6365
$ws.Range("I102").Value = 2835.2
$ws.Range("J102").Value = 1619
$ws.Range("K102").Value = 2835.2
$ws.Range("L102").Value = 1619
$ws.Range("M102").Value = -1213.2
$ws.Range("N102").Value = -4863

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 65000
$ws.Range("J64").Value = 65000
$ws.Range("L64").Value = 65000
$ws.Range("N64").Value = -65450
$ws.Range("H67").Value = 65000
$ws.Range("J67").Value = 65000
$ws.Range("L67").Value = 65000
$ws.Range("N67").Value = -66560
$ws.Range("H136").Value = 4425
$ws.Range("I136").Value = 4425
$ws.Range("K136").Value = 13275
$ws.Range("M136").Value = -10725

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2197759
$ws.Range("I96").Value = 3734584
$ws.Range("J96").Value = 2294.8572
$ws.Range("K96").Value = 3734584
$ws.Range("L96").Value = 2294.8572
$ws.Range("M96").Value = -3733211
$ws.Range("N96").Value = -5040.8572
